$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Gas6"
$ws.Range("C2").Value = "Tyro3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 20.97917533333333
$ws.Range("H2").Value = 62.93752600000001
$ws.Range("I2").Value = 0.2451892257562263
$ws.Range("J2").Value = 0.2451892257562263
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.448064
$ws.Range("N2").Value = 1.344192
$ws.Range("O2").Value = 0.113372348992564
$ws.Range("P2").Value = 0.113372348992564
$ws.Range("Q2").Value = 9.400013216554667
$ws.Range("R2").Value = 84.60011894899201
$ws.Range("S2").Value = 0.02779767847165146
$ws.Range("T2").Value = 0.02779767847165146

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Gas6"
$ws.Range("C3").Value = "Tyro3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 20.97917533333333
$ws.Range("H3").Value = 62.93752600000001
$ws.Range("I3").Value = 0.2451892257562263
$ws.Range("J3").Value = 0.2451892257562263
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.951476
$ws.Range("N3").Value = 8.854428
$ws.Range("O3").Value = 0.7468035082380574
$ws.Range("P3").Value = 0.7468035082380574
$ws.Range("Q3").Value = 61.91953249612533
$ws.Range("R3").Value = 557.2757924651281
$ws.Range("S3").Value = 0.1831081739769228
$ws.Range("T3").Value = 0.1831081739769228

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Gas6"
$ws.Range("C4").Value = "Tyro3"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 20.97917533333333
$ws.Range("H4").Value = 62.93752600000001
$ws.Range("I4").Value = 0.2451892257562263
$ws.Range("J4").Value = 0.2451892257562263
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.5526053333333333
$ws.Range("N4").Value = 1.657816
$ws.Range("O4").Value = 0.1398241427693786
$ws.Range("P4").Value = 0.1398241427693786
$ws.Range("Q4").Value = 11.59320417813511
$ws.Range("R4").Value = 104.338837603216
$ws.Range("S4").Value = 0.03428337330765197
$ws.Range("T4").Value = 0.03428337330765197

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Gas6"
$ws.Range("C5").Value = "Tyro3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 59.05285266666667
$ws.Range("H5").Value = 177.158558
$ws.Range("I5").Value = 0.6901664624076501
$ws.Range("J5").Value = 0.6901664624076501
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.448064
$ws.Range("N5").Value = 1.344192
$ws.Range("O5").Value = 0.113372348992564
$ws.Range("P5").Value = 0.113372348992564
$ws.Range("Q5").Value = 26.45945737723734
$ws.Range("R5").Value = 238.135116395136
$ws.Range("S5").Value = 0.07824579303904344
$ws.Range("T5").Value = 0.07824579303904344

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Gas6"
$ws.Range("C6").Value = "Tyro3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 59.05285266666667
$ws.Range("H6").Value = 177.158558
$ws.Range("I6").Value = 0.6901664624076501
$ws.Range("J6").Value = 0.6901664624076501
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.951476
$ws.Range("N6").Value = 8.854428
$ws.Range("O6").Value = 0.7468035082380574
$ws.Range("P6").Value = 0.7468035082380574
$ws.Range("Q6").Value = 174.2930773772027
$ws.Range("R6").Value = 1568.637696394824
$ws.Range("S6").Value = 0.5154187353942825
$ws.Range("T6").Value = 0.5154187353942825

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Gas6"
$ws.Range("C7").Value = "Tyro3"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 59.05285266666667
$ws.Range("H7").Value = 177.158558
$ws.Range("I7").Value = 0.6901664624076501
$ws.Range("J7").Value = 0.6901664624076501
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.5526053333333333
$ws.Range("N7").Value = 1.657816
$ws.Range("O7").Value = 0.1398241427693786
$ws.Range("P7").Value = 0.1398241427693786
$ws.Range("Q7").Value = 32.63292133214755
$ws.Range("R7").Value = 293.696291989328
$ws.Range("S7").Value = 0.09650193397432422
$ws.Range("T7").Value = 0.09650193397432422

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Gas6"
$ws.Range("C8").Value = "Tyro3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 5.531174333333333
$ws.Range("H8").Value = 16.593523
$ws.Range("I8").Value = 0.06464431183612354
$ws.Range("J8").Value = 0.06464431183612354
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.448064
$ws.Range("N8").Value = 1.344192
$ws.Range("O8").Value = 0.113372348992564
$ws.Range("P8").Value = 0.113372348992564
$ws.Range("Q8").Value = 2.478320096490667
$ws.Range("R8").Value = 22.304880868416
$ws.Range("S8").Value = 0.007328877481869137
$ws.Range("T8").Value = 0.007328877481869137

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Gas6"
$ws.Range("C9").Value = "Tyro3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 5.531174333333333
$ws.Range("H9").Value = 16.593523
$ws.Range("I9").Value = 0.06464431183612354
$ws.Range("J9").Value = 0.06464431183612354
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.951476
$ws.Range("N9").Value = 8.854428
$ws.Range("O9").Value = 0.7468035082380574
$ws.Range("P9").Value = 0.7468035082380574
$ws.Range("Q9").Value = 16.32512829664933
$ws.Range("R9").Value = 146.926154669844
$ws.Range("S9").Value = 0.04827659886685204
$ws.Range("T9").Value = 0.04827659886685204

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Gas6"
$ws.Range("C10").Value = "Tyro3"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 5.531174333333333
$ws.Range("H10").Value = 16.593523
$ws.Range("I10").Value = 0.06464431183612354
$ws.Range("J10").Value = 0.06464431183612354
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.5526053333333333
$ws.Range("N10").Value = 1.657816
$ws.Range("O10").Value = 0.1398241427693786
$ws.Range("P10").Value = 0.1398241427693786
$ws.Range("Q10").Value = 3.056556436196444
$ws.Range("R10").Value = 27.509007925768
$ws.Range("S10").Value = 0.009038835487402367
$ws.Range("T10").Value = 0.009038835487402367

